$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '97.963.01'
$ws.Range('E2').Value = '  +0.23%  '

$ws.Range('D3').Value = '3.363.66'
$ws.Range('E3').Value = '  -0.96%  '

$ws.Range('E4').Value = '  +0.08%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '253.01'
$ws.Range('E5').Value = '  -0.78%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '660.32'
$ws.Range('E6').Value = '  +1.47%  '

$ws.Range('E7').Value = '  -2.56%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.430'
$ws.Range('E8').Value = '  -0.09%  '

$ws.Range('E9').Value = '  +0.02%  '

$ws.Range('E10').Value = '  -3.92%  '

$ws.Range('D11').Value = '3.362.63'
$ws.Range('E11').Value = '  -0.87%  '

$ws.Range('E12').Value = '  -1.55%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '42.09'
$ws.Range('E13').Value = '  +1.56%  '

$ws.Range('D14').Value = '97.741.81'
$ws.Range('E14').Value = '  +0.42%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.14'
$ws.Range('E15').Value = '  -2.58%  '

$ws.Range('E16').Value = '  -0.81%  '

$ws.Range('D17').Value = '3.990.82'
$ws.Range('E17').Value = '  -1.23%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '8.83'
$ws.Range('E18').Value = '  +3.64%  '

$ws.Range('D19').Value = '3.359.56'
$ws.Range('E19').Value = '  -0.82%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.96'
$ws.Range('E20').Value = '  +3.28%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.538'
$ws.Range('E21').Value = '  +3.80%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.90'
$ws.Range('E22').Value = '  +1.56%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '514.80'
$ws.Range('E23').Value = '  +0.91%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.38'
$ws.Range('E24').Value = '  -1.48%  '

$ws.Range('E25').Value = '  -1.53%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.89'
$ws.Range('E26').Value = '  +11.38%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '97.10'
$ws.Range('E27').Value = '  -2.29%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.41'
$ws.Range('E28').Value = '  -2.57%  '

$ws.Range('E29').Value = '  -4.65%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '11.65'
$ws.Range('E30').Value = '  +2.22%  '

$ws.Range('E31').Value = '  +0.24%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.192'
$ws.Range('E32').Value = '  -6.52%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.60'
$ws.Range('E33').Value = '  +14.60%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('E34').Value = '  -0.05%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.572'
$ws.Range('E35').Value = '  +0.32%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '28.83'
$ws.Range('E36').Value = '  -2.55%  '

$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '8.05'
$ws.Range('E37').Value = '  +4.95%  '

$ws.Range('B38').Value = 'Fetch.AI'
$ws.Range('C38').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.53'
$ws.Range('E38').Value = '  +7.12%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '525.88'
$ws.Range('E39').Value = '  +0.77%  '

$ws.Range('E40').Value = '  -0.59%  '

$ws.Range('E41').Value = '  +0.08%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0451'
$ws.Range('E42').Value = '  +6.17%  '

$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.862'
$ws.Range('E43').Value = '  +0.80%  '

$ws.Range('B44').Value = 'WhiteBITCoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '24.44'
$ws.Range('E44').Value = '  -1.23%  '

$ws.Range('E45').Value = '  +10.33%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.72'
$ws.Range('E46').Value = '  +6.39%  '

$ws.Range('E47').Value = '  +6.26%  '

$ws.Range('E48').Value = '  -0.85%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '54.05'
$ws.Range('E49').Value = '  +5.80%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.17'

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.05'
$ws.Range('E51').Value = '  -0.69%  '
